# daily auto push: 2026-02-21 09:41 UTC
# A new sample row (2026/02/21, 土, 16, 201) was recorded and inserted into
# the data table right before the existing 2026/12/29 entries, pushing all
# subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 843; everything from the old row 843 onward shifts
# down to 844 onward (matches the diff: old row 843 -> new row 844, ...,
# old row 884 -> new row 885).
$ws.Rows("843:843").Insert()

# Column A holds date-like text (e.g. "2026/12/29") stored as plain text in
# the sheet, not as a real date serial. Force the cell to Text format first
# so Excel doesn't auto-convert the "2026/02/21" string into a date value,
# then drop back to the default "Normal" style so no stray formatting is
# left behind on the cell.
$ws.Range("A843").NumberFormat = "@"
$ws.Range("A843").Value = "2026/02/21"
$ws.Range("A843").Style = "Normal"

$ws.Range("B843").Value = "土"
$ws.Range("C843").Value = 16
$ws.Range("D843").Value = 201
